# Applies the "Initial check-in of translations changes" edit:
#  - survey!F1      "display.text"  -> "display.prompt.text"
#  - settings!C1    "display.title" -> "display.title.text"
#  - selection on settings (C15 -> C2) while leaving the sheet not active
#  - selection on survey (D39 -> F2) and make "survey" the active sheet/tab

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- Translation key renames -------------------------------------------------
# (settings' display.title.text must be appended to the shared-string table
#  before survey's display.prompt.text to match the expected string order)
$settings.Range("C1").Value = "display.title.text"
$survey.Range("F1").Value   = "display.prompt.text"

# --- Selections ---------------------------------------------------------------
# Update the settings sheet's remembered selection first (it will no longer be
# the active tab once survey is activated below).
$null = $settings.Range("C2").Select()

# Make "survey" the active sheet and set its selection/active cell.
$null = $survey.Activate()
$null = $survey.Range("F2").Select()
